$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.618.39'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').Value = '3.910.06'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''482.60'
$ws.Range('E5').Value = '  +3.61%  '
$ws.Range('D6').Value = '''147.14'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('D7').Value = '''0.617'
$ws.Range('E7').Value = '  -2.57%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -4.54%  '
$ws.Range('E10').Value = '  +7.25%  '
$ws.Range('D11').Value = '''0.0000349'
$ws.Range('E11').Value = '  +11.45%  '
$ws.Range('D12').Value = '''42.35'
$ws.Range('E12').Value = '  -3.25%  '
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').Value = '4.528.54'
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = '3.923.04'
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('D16').Value = '''14.52'
$ws.Range('E16').Value = '  -1.73%  '
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = '''19.62'
$ws.Range('E18').Value = '  -2.38%  '
$ws.Range('E19').Value = '  -3.56%  '
$ws.Range('D20').Value = '68.660.75'
$ws.Range('E20').Value = '  +1.65%  '
$ws.Range('D21').Value = '''430.72'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('D22').Value = '''14.49'
$ws.Range('E22').Value = '  -2.18%  '
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').Value = '''86.64'
$ws.Range('E24').Value = '  -2.32%  '
$ws.Range('D25').Value = '''11.40'
$ws.Range('E25').Value = '  +12.21%  '
$ws.Range('D26').Value = '''3.55'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').Value = '''10.47'
$ws.Range('E27').Value = '  +1.36%  '
$ws.Range('D28').Value = '''37.84'
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('E29').Value = '  +6.60%  '
$ws.Range('D30').Value = '''712.81'
$ws.Range('E30').Value = '  -4.16%  '
$ws.Range('E31').Value = '  -4.81%  '
$ws.Range('D32').Value = '''13.18'
$ws.Range('E32').Value = '  -4.27%  '
$ws.Range('E33').Value = '  +2.85%  '
$ws.Range('D34').Value = '0.0₃0882'
$ws.Range('E34').Value = '  +29.75%  '
$ws.Range('D35').Value = '''41.19'
$ws.Range('E35').Value = '  -4.42%  '
$ws.Range('D36').Value = '''58.25'
$ws.Range('E36').Value = '  +1.19%  '
$ws.Range('D37').Value = '''0.150'
$ws.Range('E37').Value = '  -7.55%  '
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('E40').Value = '  -2.53%  '
$ws.Range('D41').Value = '''2.75'
$ws.Range('E41').Value = '  +4.56%  '
$ws.Range('D42').Value = '''3.00'
$ws.Range('E42').Value = '  +8.58%  '
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('E44').Value = '  -3.36%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('E46').Value = '  -1.59%  '
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('E48').Value = '  +0.78%  '
$ws.Range('D49').Value = '''148.02'
$ws.Range('E49').Value = '  +2.57%  '
$ws.Range('D50').Value = '''3.12'
$ws.Range('E50').Value = '  -3.61%  '
$ws.Range('E51').Value = '  -3.36%  '
